$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price and Volume columns remain formatted as text so that
# numeric-looking strings (e.g. "2.880", "1.000") keep their exact
# textual representation instead of being coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.593.57'
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").Value = '1.882.23'
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '246.38'
$ws.Range("E5").Value = '  -0.56%  '
$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").Value = '0.4728'
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '0.2888'
$ws.Range("E8").Value = '  -1.13%  '
$ws.Range("D9").Value = '0.06535'
$ws.Range("D10").Value = '22.24'
$ws.Range("E10").Value = '  +0.99%  '
$ws.Range("D11").Value = '0.7734'
$ws.Range("E11").Value = '  +5.10%  '
$ws.Range("D12").Value = '100.89'
$ws.Range("E12").Value = '  +4.48%  '
$ws.Range("D13").Value = '0.07830'
$ws.Range("E13").Value = '  +0.24%  '
$ws.Range("D14").Value = '1.881.37'
$ws.Range("D15").Value = '5.252'
$ws.Range("E15").Value = '  +0.07%  '
$ws.Range("D16").Value = '285.87'
$ws.Range("E16").Value = '  +0.84%  '
$ws.Range("D17").Value = '30.574.25'
$ws.Range("E17").Value = '  -0.38%  '
$ws.Range("D18").Value = '13.20'
$ws.Range("E18").Value = '  -0.32%  '
$ws.Range("D19").Value = '0.000007525'
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").Value = '1.0000'
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("D21").Value = '2.126.71'
$ws.Range("E21").Value = '  -0.60%  '
$ws.Range("D22").Value = '5.367'
$ws.Range("E22").Value = '  +1.08%  '
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").Value = '6.408'
$ws.Range("E24").Value = '  +2.67%  '
$ws.Range("D25").Value = '9.132'
$ws.Range("E25").Value = '  -1.03%  '
$ws.Range("D26").Value = '163.03'
$ws.Range("E26").Value = '  -1.13%  '
$ws.Range("E27").Value = '  +0.86%  '
$ws.Range("D28").Value = '1.914'
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '1.331'
$ws.Range("E29").Value = '  -0.48%  '
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = '0.09706'
$ws.Range("E30").Value = '  -0.43%  '
$ws.Range("E31").Value = '  +1.12%  '
$ws.Range("D32").Value = '4.264'
$ws.Range("E32").Value = '  -0.79%  '
$ws.Range("D33").Value = '4.194'
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("D34").Value = '0.04849'
$ws.Range("E34").Value = '  -0.17%  '
$ws.Range("D35").Value = '1.130'
$ws.Range("E35").Value = '  +0.33%  '
$ws.Range("D36").Value = '0.6970'
$ws.Range("E36").Value = '  +0.09%  '
$ws.Range("D37").Value = '2.741'
$ws.Range("E37").Value = '  +0.62%  '
$ws.Range("E38").Value = '  +1.36%  '
$ws.Range("D39").Value = '2.880'
$ws.Range("D40").Value = '76.24'
$ws.Range("E40").Value = '  +0.33%  '
$ws.Range("D41").Value = '6.290'
$ws.Range("E41").Value = '  -1.11%  '
$ws.Range("D42").Value = '1.987'
$ws.Range("E42").Value = '  -0.56%  '
$ws.Range("D43").Value = '0.4254'
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").Value = '0.9994'
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("D45").Value = '0.8345'
$ws.Range("E45").Value = '  -0.38%  '
$ws.Range("D46").Value = '101.42'
$ws.Range("E46").Value = '  -0.14%  '
$ws.Range("D47").Value = '9.846'
$ws.Range("E47").Value = '  +3.76%  '
$ws.Range("D48").Value = '7.041'
$ws.Range("E48").Value = '  +0.23%  '
$ws.Range("D49").Value = '35.22'
$ws.Range("E49").Value = '  -1.21%  '
$ws.Range("D50").Value = '895.14'
$ws.Range("E50").Value = '  -2.30%  '
$ws.Range("D51").Value = '0.05769'
$ws.Range("E51").Value = '  +0.21%  '
